# Apply the updated cryptocurrency price/volume snapshot values.
# Source: "Updated cryptos list on Thu May 11 05:58:47 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '27.555.20'
$ws.Range('E2').Value = '  -0.62%  '

# Row 3
$ws.Range('D3').Value = '1.833.77'
$ws.Range('E3').Value = '  -0.78%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.21%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '312.71'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.38%  '

# Row 6
$ws.Range('E6').Value = '  -0.12%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4291'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.71%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3669'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +0.47%  '

# Row 9
$ws.Range('E9').Value = '  -0.76%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8652'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -1.76%  '

# Row 11
$ws.Range('E11').Value = '  -0.53%  '

# Row 12
$ws.Range('D12').Value = '1.869.61'
$ws.Range('E12').Value = '  +0.73%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.400'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +0.65%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.550'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +0.23%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.06947'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +0.16%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.003'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -0.12%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '80.72'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +1.03%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008867'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -1.72%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.001'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.06%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.39'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.15%  '

# Row 21
$ws.Range('D21').Value = '27.559.07'
$ws.Range('E21').Value = '  -0.63%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.152'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +3.31%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.85'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +4.88%  '

# Row 24
$ws.Range('D24').Value = '2.072.86'
$ws.Range('E24').Value = '  -0.48%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.992'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.10%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '154.41'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -0.86%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.89'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +1.40%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.122'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -2.47%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '114.46'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -4.35%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.830'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -2.49%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08864'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.41%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.7501'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -1.11%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.982'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.69%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.544'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.19%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.134'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +0.86%  '

# Row 36
$ws.Range('E36').Value = '  -0.06%  '

# Row 37
$ws.Range('E37').Value = '  -1.74%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05322'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -2.33%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01936'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -0.23%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.797'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -1.18%  '

# Row 41
$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.5084'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.24%  '

# Row 42
$ws.Range('B42').Value = 'Algorand'
$ws.Range('C42').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1663'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.21%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.539'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -1.55%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.327'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.85%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.45'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +0.14%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '105.76'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -0.07%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.06485'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -1.06%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4684'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.48%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.9999'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.12%  '

# Row 50
$ws.Range('E50').Value = '  -1.73%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '63.69'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -1.59%  '
